$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WR")

# Add the new player row (Week 16 logging added a new player: E.Wolf)
$ws.Range("A12").Value = "E.Wolf"
$ws.Range("B12:J12").Value = 0

# Move the selection to where the user's cursor ended up (K13)
$ws.Range("K13").Select()
